# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml  "Office Theme" (clrScheme "Office")      -> becomes the "Integral" / "Red Violet" theme
#   ppt/theme/theme2.xml  "Integral"     (clrScheme "Red Violet")  -> becomes the "Office Theme" / "Office" theme
#
# theme2.xml is the theme that actually drives the deck (it's the one referenced
# by the slide master / presentation.xml.rels), so that's the theme PowerPoint's
# object model exposes as $p.SlideMaster.Theme (Master.Theme). Its fontScheme and
# fmtScheme (gradients/lines/effects) are already byte-identical between the two
# theme parts - only the 12 clrScheme colors differ - so re-colouring that theme
# with the "Office Theme" palette reproduces the target theme2.xml content.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

function ConvertHexToVbRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme colour scheme, in dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order
# (ThemeColorScheme.Item(1..12)).
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertHexToVbRgb $officeThemeColors[$i - 1]
}
